$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Anxa2"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 105.2534506666667
$ws.Range("H2").Value = 315.760352
$ws.Range("I2").Value = 0.2966477300323703
$ws.Range("J2").Value = 0.2966477300323703
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 51.23401333333334
$ws.Range("N2").Value = 153.70204
$ws.Range("O2").Value = 0.9688226788583661
$ws.Range("P2").Value = 0.9688226788583661
$ws.Range("Q2").Value = 5392.556694835343
$ws.Range("R2").Value = 48533.01025351808
$ws.Range("S2").Value = 0.2873990484872144
$ws.Range("T2").Value = 0.2873990484872144

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Anxa2"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 105.2534506666667
$ws.Range("H3").Value = 315.760352
$ws.Range("I3").Value = 0.2966477300323703
$ws.Range("J3").Value = 0.2966477300323703
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.648742666666666
$ws.Range("N3").Value = 4.946228
$ws.Range("O3").Value = 0.03117732114163389
$ws.Range("P3").Value = 0.03117732114163389
$ws.Range("Q3").Value = 173.5358549280284
$ws.Range("R3").Value = 1561.822694352256
$ws.Range("S3").Value = 0.00924868154515592
$ws.Range("T3").Value = 0.009248681545155922

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Anxa2"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 161.0956266666667
$ws.Range("H4").Value = 483.28688
$ws.Range("I4").Value = 0.4540340641196984
$ws.Range("J4").Value = 0.4540340641196984
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 51.23401333333334
$ws.Range("N4").Value = 153.70204
$ws.Range("O4").Value = 0.9688226788583661
$ws.Range("P4").Value = 0.9688226788583661
$ws.Range("Q4").Value = 8253.57548458169
$ws.Range("R4").Value = 74282.1793612352
$ws.Range("S4").Value = 0.4398784982933974
$ws.Range("T4").Value = 0.4398784982933974

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Anxa2"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 161.0956266666667
$ws.Range("H5").Value = 483.28688
$ws.Range("I5").Value = 0.4540340641196984
$ws.Range("J5").Value = 0.4540340641196984
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.648742666666666
$ws.Range("N5").Value = 4.946228
$ws.Range("O5").Value = 0.03117732114163389
$ws.Range("P5").Value = 0.03117732114163389
$ws.Range("Q5").Value = 265.6052330987378
$ws.Range("R5").Value = 2390.44709788864
$ws.Range("S5").Value = 0.01415556582630103
$ws.Range("T5").Value = 0.01415556582630103

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Anxa2"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 88.46048300000001
$ws.Range("H6").Value = 265.381449
$ws.Range("I6").Value = 0.2493182058479313
$ws.Range("J6").Value = 0.2493182058479313
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 51.23401333333334
$ws.Range("N6").Value = 153.70204
$ws.Range("O6").Value = 0.9688226788583661
$ws.Range("P6").Value = 0.9688226788583661
$ws.Range("Q6").Value = 4532.185565495108
$ws.Range("R6").Value = 40789.67008945597
$ws.Range("S6").Value = 0.2415451320777544
$ws.Range("T6").Value = 0.2415451320777544

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Anxa2"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 88.46048300000001
$ws.Range("H7").Value = 265.381449
$ws.Range("I7").Value = 0.2493182058479313
$ws.Range("J7").Value = 0.2493182058479313
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 1.648742666666666
$ws.Range("N7").Value = 4.946228
$ws.Range("O7").Value = 0.03117732114163389
$ws.Range("P7").Value = 0.03117732114163389
$ws.Range("Q7").Value = 145.8485726360413
$ws.Range("R7").Value = 1312.637153724372
$ws.Range("S7").Value = 0.007773073770176938
$ws.Range("T7").Value = 0.007773073770176939

Write-Output "done"